# Updated legacy GSC export data:
# - The oldest day (2025-10-14) is dropped from the Coverage/"Chart" sheet
#   and every later day's row shifts up by one.
# - The three most-recent days (now rows 2-4) don't have "Not indexed" /
#   "Indexed" counts yet, so those two columns are cleared to blank for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the first data row (2025-10-14); everything below shifts up one row.
$ws.Rows.Item(2).Delete()

# The three newest rows (2025-10-15, 2025-10-16, 2025-10-17) have no
# "Not indexed" / "Indexed" figures yet - blank them out.
$ws.Range("B2:C4").Value = ""
